$wb = $excel.ActiveWorkbook

# "Generate Report for Archive": the tracked item has moved on from
# "Ready for handoff" to "In Translation". Update the status text
# everywhere it is shown -- the Overview roll-up (one column per locale)
# and each locale sheet's own "Status" column -- then re-fit the columns
# that used to hold the longer "Ready for handoff" text.

$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "In Translation"
$ov.Range("F2").Value = "In Translation"

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "In Translation"

$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "In Translation"

$ov.Range("E1").ColumnWidth = 12.5
$ov.Range("F1").ColumnWidth = 12.5
$zh.Range("C1").ColumnWidth = 12.5
$de.Range("C1").ColumnWidth = 12.5
